$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Content")

# ---------------------------------------------------------------------
# Row 5: existing "JLoBeauty" row becomes "DrDenese" (new shared string).
# ---------------------------------------------------------------------
$ws.Range("B5").Value = "DrDenese"

# ---------------------------------------------------------------------
# Three brand-new rows (6,7,8) are inserted ahead of the old row-6 "End"
# marker, which is re-created at row 9. Format the new rows to match the
# existing data rows (fill style index 2) before writing values so the
# shared-string table gets populated in column-major order: all of
# column B top-to-bottom, then column C, matching the authored file
# (DrDenese, Smileactives, MallyBeauty, then Core2).
# ---------------------------------------------------------------------
$ws.Range("A2:C2").Copy()
$ws.Range("A6:C8").PasteSpecial(-4122)
$ws.Range("A6").Copy()
$ws.Range("A9").PasteSpecial(-4122)

$ws.Range("B6").Value = "Smileactives"
$ws.Range("B7").Value = "MallyBeauty"
$ws.Range("B8").Value = "JLoBeauty"

$ws.Range("C6").Value = "Core2"
$ws.Range("C7").Value = "Core"
$ws.Range("C8").Value = "Core"

$ws.Range("A6").Value = "Prod"
$ws.Range("A7").Value = "Prod"
$ws.Range("A8").Value = "Prod"

$ws.Range("A9").Value = "End"

# ---------------------------------------------------------------------
# Old block (rows 15-19) collapses: row 15 disappears, rows 16-17 become
# JLoBeauty/Core + End (with only column A populated), rows 18-19
# disappear. Delete the whole contiguous block and re-insert blank rows
# in its place so the removed rows leave no trace in the sheet, then
# rebuild just the two surviving rows.
# ---------------------------------------------------------------------
$ws.Rows("15:19").Delete()
$ws.Rows("15:19").Insert()

$ws.Range("A2:C2").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)
$ws.Range("A6").Copy()
$ws.Range("A17").PasteSpecial(-4122)

$ws.Range("A16").Value = "Prod"
$ws.Range("B16").Value = "JLoBeauty"
$ws.Range("C16").Value = "Core"

$ws.Range("A17").Value = "End"

# ---------------------------------------------------------------------
# New block, rows 21-23.
# ---------------------------------------------------------------------
$ws.Range("A2:C2").Copy()
$ws.Range("A21:C23").PasteSpecial(-4122)

$ws.Range("A21").Value = "Prod"
$ws.Range("B21").Value = "Smileactives"
$ws.Range("C21").Value = "Core"

$ws.Range("A22").Value = "Prod"
$ws.Range("B22").Value = "CrepeErase"
$ws.Range("C22").Value = "Core"

$ws.Range("A23").Value = "Prod"
$ws.Range("B23").Value = "MeaningfulBeauty"
$ws.Range("C23").Value = "core_full_30_90"

# ---------------------------------------------------------------------
# New block, rows 28-35.
# ---------------------------------------------------------------------
$ws.Range("A2:C2").Copy()
$ws.Range("A28:C34").PasteSpecial(-4122)
$ws.Range("A6").Copy()
$ws.Range("A35").PasteSpecial(-4122)

$ws.Range("A28").Value = "Prod"
$ws.Range("B28").Value = "CrepeErase"
$ws.Range("C28").Value = "Core"

$ws.Range("A29").Value = "Prod"
$ws.Range("B29").Value = "MeaningfulBeauty"
$ws.Range("C29").Value = "Core"

$ws.Range("A30").Value = "Prod"
$ws.Range("B30").Value = "WestmoreBeauty"
$ws.Range("C30").Value = "Core"

$ws.Range("A31").Value = "Prod"
$ws.Range("B31").Value = "DrDenese"
$ws.Range("C31").Value = "Core"

$ws.Range("A32").Value = "Prod"
$ws.Range("B32").Value = "Smileactives"
$ws.Range("C32").Value = "Core2"

$ws.Range("A33").Value = "Prod"
$ws.Range("B33").Value = "MallyBeauty"
$ws.Range("C33").Value = "Core"

$ws.Range("A34").Value = "Prod"
$ws.Range("B34").Value = "JLoBeauty"
$ws.Range("C34").Value = "Core"

$ws.Range("A35").Value = "End"

# Match the final selection recorded in the authored workbook.
$ws.Range("A2:XFD15").Select()
